# rebase one of the demo measures
# - Update the attendances report_comment text (E2) to mention the re-base.
# - Widen column E to fit the longer comment text.
# - Clear the lingering cell selection (E3) left over in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "This is a comment about the attendances metric, which has been re-based as a demonstration.  This text is added via 'report_config.xlsx'"

# Column E's stored <col> width needs to land on 114 after Excel's internal
# padding (+5/6 character) is applied by the ColumnWidth setter.
$ws.Columns("E").ColumnWidth = 113.16666666666667

# Reset the active selection to A1 so the sheet no longer reports the old
# E3 selection.
$ws.Range("A1").Select() | Out-Null
